# Updating attendance files 4th march
# Fill in "Session 2" (column F) attendance marks for rows 7-80
# (S.No. 1-74). Column J ("Total Present") auto-recalculates via its
# existing SUM(E:I) formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Most participants were present for session 2 (value 1); fill the whole
# block first, then patch the few absentees below.
$ws.Range("F7:F80").Value = 1

# Absentees for session 2
$ws.Range("F20").Value = 0
$ws.Range("F52").Value = 0
$ws.Range("F73").Value = 0
